$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55; this shifts existing rows 55-80 down to 56-81
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with the new data record
$ws.Range("A55").Value = 9
$ws.Range("B55").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = 44523
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 100112022
$ws.Range("G55").Value = "Arveja Verde"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 34
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 16000
$ws.Range("M55").Value = 15500
$ws.Range("N55").Value = '$/saco 25 kilos'
$ws.Range("O55").Value = "Región Metropolitana"
$ws.Range("P55").Value = 620
$ws.Range("Q55").Value = 25
$ws.Range("R55").Value = "Hortaliza"
